$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format for this update so that numeric-looking
# values (e.g. "1.001", "306.96") are stored as text, matching the
# original inlineStr/text cell content instead of being auto-converted
# to numbers by Excel smart input parsing.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.152.20"

$ws.Range("D3").Value = "1.871.52"
$ws.Range("E3").Value = "  -1.90%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "306.96"
$ws.Range("E5").Value = "  -1.81%  "

$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").Value = "0.5126"
$ws.Range("E7").Value = "  +2.61%  "

$ws.Range("D8").Value = "0.3748"
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").Value = "0.07134"
$ws.Range("E9").Value = "  -1.89%  "

$ws.Range("D10").Value = "0.8875"
$ws.Range("E10").Value = "  -2.34%  "

$ws.Range("D11").Value = "20.65"
$ws.Range("E11").Value = "  -2.77%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07540"
$ws.Range("E12").Value = "  -1.32%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.857.08"
$ws.Range("E13").Value = "  -2.83%  "

$ws.Range("D14").Value = "5.329"
$ws.Range("E14").Value = "  -2.50%  "

$ws.Range("D15").Value = "89.18"
$ws.Range("E15").Value = "  -3.40%  "

$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "0.000008477"
$ws.Range("E17").Value = "  -2.76%  "

$ws.Range("E18").Value = "  -3.83%  "

$ws.Range("D19").Value = "0.9994"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "27.220.23"
$ws.Range("E20").Value = "  -2.45%  "

$ws.Range("D21").Value = "5.054"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("D22").Value = "2.109.06"
$ws.Range("E22").Value = "  -2.99%  "

$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("D24").Value = "6.474"
$ws.Range("E24").Value = "  -1.96%  "

$ws.Range("D25").Value = "149.96"
$ws.Range("E25").Value = "  -1.72%  "

$ws.Range("D26").Value = "1.846"
$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").Value = "17.96"
$ws.Range("E27").Value = "  -2.27%  "

$ws.Range("D28").Value = "2.097"
$ws.Range("E28").Value = "  -5.50%  "

$ws.Range("D29").Value = "112.99"
$ws.Range("E29").Value = "  -1.80%  "

$ws.Range("D30").Value = "4.722"
$ws.Range("E30").Value = "  -3.27%  "

$ws.Range("D31").Value = "4.670"
$ws.Range("E31").Value = "  -2.95%  "

$ws.Range("D32").Value = "0.09025"
$ws.Range("E32").Value = "  +0.69%  "

$ws.Range("D33").Value = "0.05133"
$ws.Range("E33").Value = "  -2.76%  "

$ws.Range("D34").Value = "3.091"
$ws.Range("E34").Value = "  -3.09%  "

$ws.Range("D35").Value = "1.157"
$ws.Range("E35").Value = "  -6.38%  "

$ws.Range("D36").Value = "0.7362"
$ws.Range("E36").Value = "  -6.51%  "

$ws.Range("D37").Value = "0.02056"
$ws.Range("E37").Value = "  -1.08%  "

$ws.Range("D38").Value = "2.517"
$ws.Range("E38").Value = "  -5.18%  "

$ws.Range("D39").Value = "3.066"
$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("D40").Value = "1.073"
$ws.Range("E40").Value = "  -1.62%  "

$ws.Range("D41").Value = "0.5365"
$ws.Range("E41").Value = "  -2.89%  "

$ws.Range("D42").Value = "6.586"
$ws.Range("E42").Value = "  -3.22%  "

$ws.Range("D43").Value = "117.26"
$ws.Range("E43").Value = "  +3.05%  "

$ws.Range("D44").Value = "8.332"
$ws.Range("E44").Value = "  -2.05%  "

$ws.Range("E45").Value = "  -2.40%  "

$ws.Range("D46").Value = "0.4639"
$ws.Range("E46").Value = "  -3.68%  "

$ws.Range("D47").Value = "0.9988"
$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").Value = "10.09"
$ws.Range("E48").Value = "  -4.78%  "

$ws.Range("D49").Value = "1.570"
$ws.Range("E49").Value = "  -3.93%  "

$ws.Range("D50").Value = "64.39"
$ws.Range("E50").Value = "  -4.31%  "

$ws.Range("D51").Value = "36.53"
$ws.Range("E51").Value = "  -1.01%  "

# Restore the original (default) cell style on column D now that the
# values have been entered as text, so no stray number format remains
# applied to the cells.
$dRange.Style = "Normal"
